$d = $word.ActiveDocument

# Replacement pairs (old -> new). Order matters: 39÷4= -> 70÷2= must run
# before 76÷8= -> 39÷4= so the freshly written "39÷4=" text from the
# latter isn't re-matched by the former.
$pairs = @(
    @("86÷9=", "75÷3="),
    @("60÷2=", "13÷3="),
    @("17÷6=", "63÷8="),
    @("75÷6=", "33÷2="),
    @("39÷4=", "70÷2="),
    @("76÷8=", "39÷4="),
    @("45÷9=", "73÷2="),
    @("92÷7=", "76÷5="),
    @("44÷4=", "83÷3="),
    @("19÷4=", "89÷8="),
    @("93÷7=", "43÷4="),
    @("18÷5=", "75÷9="),
    @("47÷4=", "63÷6="),
    @("37÷8=", "22÷7="),
    @("61÷9=", "60÷6="),
    @("31÷3=", "87÷8="),
    @("82÷3=", "79÷4="),
    @("66÷5=", "15÷7="),
    @("12÷2=", "37÷2="),
    @("95÷8=", "94÷7="),
    @("60÷4=", "88÷2="),
    @("94÷8=", "98÷9="),
    @("48÷8=", "68÷4="),
    @("32÷7=", "99÷3="),
    @("93÷2=", "79÷2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
